$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week date range) ---
# Only the trailing numeric run in each shared string changes; use Characters()
# to edit just that substring in place and keep the rest of the text intact.
$ws.Range("A8").Characters(21, 2).Text = "37"

$ws.Range("C9").Characters(27, 8).Text = "9/8/2025"
$ws.Range("C9").Characters(46, 8).Text = "9/14/2025"

# --- Crime-complaint numeric table updates (rows 14-30, 33) ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = -85.714285714285
$ws.Range("F14").Value = 13
$ws.Range("G14").Value = 13
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 79
$ws.Range("J14").Value = 88
$ws.Range("K14").Value = -10.227272727272
$ws.Range("L14").Value = -21
$ws.Range("M14").Value = -18.556701030927
$ws.Range("N14").Value = -77.620396600566
$ws.Range("C15").Value = 13
$ws.Range("D15").Value = 12
$ws.Range("E15").Value = 8.333333333333
$ws.Range("F15").Value = 46
$ws.Range("G15").Value = 35
$ws.Range("H15").Value = 31.428571428571
$ws.Range("I15").Value = 380
$ws.Range("J15").Value = 298
$ws.Range("K15").Value = 27.516778523489
$ws.Range("L15").Value = 37.184115523465
$ws.Range("M15").Value = 71.171171171171
$ws.Range("N15").Value = -25.78125
$ws.Range("C16").Value = 101
$ws.Range("D16").Value = 98
$ws.Range("E16").Value = 3.061224489795
$ws.Range("G16").Value = 422
$ws.Range("H16").Value = -14.928909952606
$ws.Range("I16").Value = 3411
$ws.Range("J16").Value = 3532
$ws.Range("K16").Value = -3.425821064552
$ws.Range("L16").Value = -0.553935860058
$ws.Range("M16").Value = 12.947019867549
$ws.Range("N16").Value = -69.663820704375
$ws.Range("C17").Value = 162
$ws.Range("D17").Value = 182
$ws.Range("E17").Value = -10.989010989011
$ws.Range("G17").Value = 698
$ws.Range("H17").Value = -4.154727793696
$ws.Range("I17").Value = 6225
$ws.Range("J17").Value = 5989
$ws.Range("K17").Value = 3.940557689096
$ws.Range("L17").Value = 6.446648426812
$ws.Range("M17").Value = 96.558257025576
$ws.Range("N17").Value = -5.150083803138
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 53
$ws.Range("E18").Value = -11.320754716981
$ws.Range("F18").Value = 192
$ws.Range("G18").Value = 220
$ws.Range("H18").Value = -12.727272727272
$ws.Range("I18").Value = 1984
$ws.Range("J18").Value = 2075
$ws.Range("K18").Value = -4.385542168674
$ws.Range("L18").Value = -5.837683910773
$ws.Range("M18").Value = -12.444836716681
$ws.Range("N18").Value = -85.016237444301
$ws.Range("C19").Value = 194
$ws.Range("D19").Value = 183
$ws.Range("E19").Value = 6.010928961748
$ws.Range("F19").Value = 743
$ws.Range("G19").Value = 787
$ws.Range("H19").Value = -5.590851334180
$ws.Range("I19").Value = 6479
$ws.Range("J19").Value = 6606
$ws.Range("K19").Value = -1.922494701786
$ws.Range("L19").Value = 14.713172804532
$ws.Range("M19").Value = 99.047619047619
$ws.Range("N19").Value = 21.306871372402
$ws.Range("C20").Value = 61
$ws.Range("D20").Value = 86
$ws.Range("E20").Value = -29.069767441860
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 363
$ws.Range("H20").Value = -18.457300275482
$ws.Range("I20").Value = 3178
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 5.933333333333
$ws.Range("L20").Value = -15.970386039132
$ws.Range("M20").Value = 117.820424948595
$ws.Range("N20").Value = -70.428956918209
$ws.Range("C21").Value = 579
$ws.Range("D21").Value = 621
$ws.Range("E21").Value = -6.763285024154
$ws.Range("F21").Value = 2318
$ws.Range("G21").Value = 2538
$ws.Range("H21").Value = -8.668242710795
$ws.Range("I21").Value = 21736
$ws.Range("J21").Value = 21588
$ws.Range("K21").Value = 0.685566055215
$ws.Range("L21").Value = 2.567006417516
$ws.Range("M21").Value = 61.174551386623
$ws.Range("N21").Value = -54.717610049790
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -20
$ws.Range("F22").Value = 25
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 217
$ws.Range("J22").Value = 240
$ws.Range("K22").Value = -9.583333333333
$ws.Range("L22").Value = 4.326923076923
$ws.Range("M22").Value = -0.913242009132
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 28
$ws.Range("E23").Value = 10.714285714285
$ws.Range("F23").Value = 103
$ws.Range("G23").Value = 154
$ws.Range("H23").Value = -33.116883116883
$ws.Range("I23").Value = 1104
$ws.Range("J23").Value = 1246
$ws.Range("K23").Value = -11.396468699839
$ws.Range("L23").Value = -13.75
$ws.Range("M23").Value = 46.419098143236
$ws.Range("C24").Value = 326
$ws.Range("D24").Value = 336
$ws.Range("E24").Value = -2.976190476190
$ws.Range("F24").Value = 1344
$ws.Range("G24").Value = 1299
$ws.Range("H24").Value = 3.464203233256
$ws.Range("I24").Value = 12635
$ws.Range("J24").Value = 11440
$ws.Range("K24").Value = 10.445804195804
$ws.Range("L24").Value = -0.307716585134
$ws.Range("M24").Value = 40.905542544886
$ws.Range("C25").Value = 92
$ws.Range("D25").Value = 124
$ws.Range("E25").Value = -25.806451612903
$ws.Range("F25").Value = 443
$ws.Range("G25").Value = 494
$ws.Range("H25").Value = -10.323886639676
$ws.Range("I25").Value = 4173
$ws.Range("J25").Value = 4553
$ws.Range("K25").Value = -8.346145398638
$ws.Range("L25").Value = -23.078341013824
$ws.Range("C26").Value = 204
$ws.Range("D26").Value = 225
$ws.Range("E26").Value = -9.333333333333
$ws.Range("F26").Value = 879
$ws.Range("G26").Value = 853
$ws.Range("H26").Value = 3.048065650644
$ws.Range("I26").Value = 7911
$ws.Range("J26").Value = 7879
$ws.Range("K26").Value = 0.406142911536
$ws.Range("L26").Value = 6.287787182587
$ws.Range("M26").Value = 0.139240506329
$ws.Range("C27").Value = 21
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = 31.25
$ws.Range("F27").Value = 64
$ws.Range("G27").Value = 45
$ws.Range("H27").Value = 42.222222222222
$ws.Range("I27").Value = 484
$ws.Range("J27").Value = 457
$ws.Range("K27").Value = 5.908096280087
$ws.Range("L27").Value = 3.862660944206
$ws.Range("C28").Value = 27
$ws.Range("D28").Value = 23
$ws.Range("E28").Value = 17.391304347826
$ws.Range("F28").Value = 84
$ws.Range("G28").Value = 83
$ws.Range("H28").Value = 1.204819277108
$ws.Range("I28").Value = 768
$ws.Range("J28").Value = 849
$ws.Range("K28").Value = -9.540636042402
$ws.Range("L28").Value = 4.489795918367
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 17
$ws.Range("E29").Value = -82.352941176470
$ws.Range("F29").Value = 40
$ws.Range("G29").Value = 46
$ws.Range("H29").Value = -13.043478260869
$ws.Range("I29").Value = 242
$ws.Range("J29").Value = 314
$ws.Range("K29").Value = -22.929936305732
$ws.Range("L29").Value = -17.966101694915
$ws.Range("M29").Value = -30.857142857142
$ws.Range("N29").Value = -75.992063492063
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 14
$ws.Range("E30").Value = -85.714285714285
$ws.Range("F30").Value = 22
$ws.Range("G30").Value = 36
$ws.Range("H30").Value = -38.888888888888
$ws.Range("I30").Value = 195
$ws.Range("J30").Value = 247
$ws.Range("K30").Value = -21.052631578947
$ws.Range("L30").Value = -19.421487603305
$ws.Range("M30").Value = -33.219178082191
$ws.Range("N30").Value = -78.711790393013
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = -40
$ws.Range("I33").Value = 22
$ws.Range("K33").Value = -42.105263157894
$ws.Range("L33").Value = -37.142857142857

# --- C33 special case: value becomes literal text "0" (shared string),
#     matching the style already used by other text cells in this table (e.g. D31). ---
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "0"
$ws.Range("D31").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$excel.CutCopyMode = 0
